$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 128.125
$ws.Range("I2").Value = 118.1
$ws.Range("K2").Value = 118.1
$ws.Range("M2").Value = -5.099999999999994
$ws.Range("H3").Value = 24758.143
$ws.Range("J3").Value = 24758.143
$ws.Range("L3").Value = 24758.143
$ws.Range("N3").Value = -24986.143
$ws.Range("H17").Value = 2101.0322
$ws.Range("J17").Value = 2389.28
$ws.Range("L17").Value = 7167.84
$ws.Range("N17").Value = -7503.84
$ws.Range("H33").Value = 73
$ws.Range("I33").Value = 73
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 73
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 156
$ws.Range("H51").Value = 2984.7
$ws.Range("J51").Value = 3066.611
$ws.Range("L51").Value = 3066.611
$ws.Range("N51").Value = -4034.611
$ws.Range("H86").Value = 1280.3
$ws.Range("I86").Value = 959.8
$ws.Range("K86").Value = 959.8
$ws.Range("M86").Value = 163.2
$ws.Range("H89").Value = 1280.3
$ws.Range("I89").Value = 959.8
$ws.Range("K89").Value = 4799
$ws.Range("M89").Value = 817
$ws.Range("H96").Value = 111111310
$ws.Range("I96").Value = 111111310
$ws.Range("K96").Value = 333333930
$ws.Range("M96").Value = -333332557
$ws.Range("H98").Value = 371.44446
$ws.Range("I98").Value = 371.44446
$ws.Range("K98").Value = 371.44446
$ws.Range("M98").Value = 1126.55554
$ws.Range("H102").Value = 24758.143
$ws.Range("J102").Value = 24758.143
$ws.Range("L102").Value = 24758.143
$ws.Range("N102").Value = -31248.143
$ws.Range("H104").Value = 977.6
$ws.Range("I104").Value = 972
$ws.Range("K104").Value = 2916
$ws.Range("M104").Value = -1169
$ws.Range("H105").Value = 24264.834
$ws.Range("J105").Value = 24264.834
$ws.Range("L105").Value = 24264.834
$ws.Range("N105").Value = -31252.834
$ws.Range("H113").Value = 1736
$ws.Range("J113").Value = 1495
$ws.Range("L113").Value = 1495
$ws.Range("H122").Value = 371.44446
$ws.Range("I122").Value = 371.44446
$ws.Range("K122").Value = 1114.33338
$ws.Range("M122").Value = 1335.66662
$ws.Range("H132").Value = 10901.762
$ws.Range("I132").Value = 10901.762
$ws.Range("K132").Value = 32705.286
$ws.Range("M132").Value = -30175.286
$ws.Range("H138").Value = 6430.3335
$ws.Range("I138").Value = 499.33334
$ws.Range("J138").Value = 8407.333000000001
$ws.Range("K138").Value = 1498.00002
$ws.Range("L138").Value = 25221.999
$ws.Range("M138").Value = 3641.99998
$ws.Range("N138").Value = -35501.999
$ws.Range("N113").Value = -8003
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5996.5
$ws.Range("I61").Value = 3989.5
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 3989.5
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -3777.5
$ws.Range("N61").Value = -7424
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("H132").Value = 1475.625
$ws.Range("I132").Value = 1460.3846
$ws.Range("J132").Value = 1541.6666
$ws.Range("K132").Value = 4381.1538
$ws.Range("L132").Value = 4624.9998
$ws.Range("M132").Value = -1851.1538
$ws.Range("N132").Value = -9684.9998
$ws.Range("H136").Value = 5996.5
$ws.Range("I136").Value = 3989.5
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 11968.5
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -9418.5
$ws.Range("N136").Value = -26100
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 56.96154
$ws.Range("I7").Value = 32.882355
$ws.Range("J7").Value = 102.44444
$ws.Range("K7").Value = 32.882355
$ws.Range("L7").Value = 102.44444
$ws.Range("M7").Value = 80.11764500000001
$ws.Range("N7").Value = -328.44444
$ws.Range("H31").Value = 5557.516
$ws.Range("I31").Value = 2765.7222
$ws.Range("K31").Value = 2765.7222
$ws.Range("M31").Value = -2470.7222
$ws.Range("H34").Value = 5557.516
$ws.Range("I34").Value = 2765.7222
$ws.Range("K34").Value = 2765.7222
$ws.Range("M34").Value = -2563.7222
$ws.Range("H106").Value = 38199.75
$ws.Range("J106").Value = 38199.75
$ws.Range("L106").Value = 38199.75
$ws.Range("N106").Value = -40723.75
$ws.Range("H107").Value = 496.33334
$ws.Range("I107").Value = 496.33334
$ws.Range("K107").Value = 496.33334
$ws.Range("M107").Value = 1423.66666
$ws.Range("H132").Value = 1306.6129
$ws.Range("I132").Value = 1306.6129
$ws.Range("K132").Value = 3919.8387
$ws.Range("M132").Value = -1389.8387
$ws.Range("H138").Value = 105000
$ws.Range("J138").Value = 105000
$ws.Range("L138").Value = 105000
$ws.Range("N138").Value = -115280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 325
$ws.Range("J23").Value = 450
$ws.Range("L23").Value = 450
$ws.Range("N23").Value = -896
$ws.Range("H80").Value = 2824
$ws.Range("I80").Value = 2711.125
$ws.Range("J80").Value = 3049.75
$ws.Range("K80").Value = 2711.125
$ws.Range("L80").Value = 3049.75
$ws.Range("M80").Value = -1713.125
$ws.Range("N80").Value = -5045.75
$ws.Range("H83").Value = 2824
$ws.Range("I83").Value = 2711.125
$ws.Range("J83").Value = 3049.75
$ws.Range("K83").Value = 13555.625
$ws.Range("L83").Value = 15248.75
$ws.Range("M83").Value = -8563.625
$ws.Range("N83").Value = -25232.75
$ws.Range("H102").Value = 2353.875
$ws.Range("I102").Value = 1462.2
$ws.Range("K102").Value = 1462.2
$ws.Range("M102").Value = 159.8
$ws.Range("H132").Value = 57761.723
$ws.Range("I132").Value = 57761.723
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 173285.169
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -170755.169
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3890.5
$ws.Range("I7").Value = 2210.75
$ws.Range("K7").Value = 2210.75
$ws.Range("M7").Value = -2098.75
$ws.Range("H16").Value = 2950
$ws.Range("I16").Value = 2950
$ws.Range("K16").Value = 2950
$ws.Range("M16").Value = -2780
$ws.Range("H40").Value = 4620.8335
$ws.Range("I40").Value = 4422.0586
$ws.Range("K40").Value = 4422.0586
$ws.Range("M40").Value = -4286.0586
$ws.Range("H46").Value = 6799.8335
$ws.Range("I46").Value = 799
$ws.Range("J46").Value = 8000
$ws.Range("K46").Value = 799
$ws.Range("L46").Value = 8000
$ws.Range("M46").Value = -611
$ws.Range("N46").Value = -8376
$ws.Range("H93").Value = 1445.6875
$ws.Range("I93").Value = 1410.9166
$ws.Range("J93").Value = 1550
$ws.Range("K93").Value = 1410.9166
$ws.Range("L93").Value = 1550
$ws.Range("M93").Value = -162.9166
$ws.Range("N93").Value = -4046
$ws.Range("H100").Value = 10000
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("H126").Value = 3890.5
$ws.Range("I126").Value = 2210.75
$ws.Range("K126").Value = 6632.25
$ws.Range("M126").Value = -4162.25
$ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5605
$ws.Range("I122").Value = 4635
$ws.Range("J122").Value = 6575
$ws.Range("K122").Value = 13905
$ws.Range("L122").Value = 19725
$ws.Range("M122").Value = -11455
$ws.Range("N122").Value = -24625
$ws.Range("H126").Value = 2955.6843
$ws.Range("I126").Value = 1212.9231
$ws.Range("K126").Value = 3638.7693
$ws.Range("M126").Value = -1168.7693
